# New crime data collected - weekly CompStat update for 123rd Precinct
# Updates: Volume/Number header, report week dates, and weekly crime
# complaint stats (rows 15-28, 33 of the "Crime Complaints" table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header rich text: "Volume 32   Number  35" -> "...  36"
#    and "Report Covering the Week  8/25/2025  Through  8/31/2025"
#    -> "...  9/1/2025  Through  9/7/2025"
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "36"

# Replace the second date first so the first date's character offsets
# (which come before it) stay valid.
$ws.Range("C9").Characters(47, 9).Text = "9/7/2025"
$ws.Range("C9").Characters(27, 9).Text = "9/1/2025"

# ---------------------------------------------------------------------
# Helper donors (untouched row 14) used to restore number formats when
# a cell's underlying type flips between text ("0" / "***.*") and a
# real number - assigning .Value alone does not change the style, and
# assigning a text value that looks numeric gets auto-typed as a
# number, so we force text via NumberFormat, set the value, then copy
# the real number format back in with PasteSpecial.
# ---------------------------------------------------------------------
function Set-TextCell($cell, $text, $donor) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $donor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

function Set-NumberCell($cell, $number, $donor) {
    $cell.Value = $number
    $donor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

$donorText0 = $ws.Range("C14")      # style 13, text "0"
$donorTextStar = $ws.Range("E14")   # style 13, text "***.*"
$donorInt = $ws.Range("J14")        # style 14, integer
$donorPct = $ws.Range("K14")        # style 15, percent

# ---------------------------------------------------------------------
# Row 15 (Murder): G15 1 -> "0", H15 -100 -> "***.*"
# ---------------------------------------------------------------------
Set-TextCell $ws.Range("G15") "0" $donorText0
Set-TextCell $ws.Range("H15") "***.*" $donorTextStar

# ---------------------------------------------------------------------
# Row 16 (Rape): C16 1 -> "0", F16 2 -> 1
# ---------------------------------------------------------------------
Set-TextCell $ws.Range("C16") "0" $donorText0
$ws.Range("F16").Value = 1
$ws.Range("M16").Value = 5.555555555555
$ws.Range("N16").Value = -53.658536585365

# ---------------------------------------------------------------------
# Row 17 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 74
$ws.Range("J17").Value = 65
$ws.Range("K17").Value = 13.846153846153
$ws.Range("L17").Value = 15.625
$ws.Range("M17").Value = 131.25
$ws.Range("N17").Value = -6.329113924050

# ---------------------------------------------------------------------
# Row 18 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 40
$ws.Range("L18").Value = -10.256410256410
$ws.Range("M18").Value = -52.054794520547
$ws.Range("N18").Value = -83.253588516746

# ---------------------------------------------------------------------
# Row 19 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 32
$ws.Range("H19").Value = 52.380952380952
$ws.Range("I19").Value = 204
$ws.Range("J19").Value = 203
$ws.Range("K19").Value = 0.492610837438
$ws.Range("L19").Value = 4.081632653061
$ws.Range("M19").Value = 110.309278350515
$ws.Range("N19").Value = 32.467532467532

# ---------------------------------------------------------------------
# Row 20 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("D20").Value = 1
$ws.Range("G20").Value = 4
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -41.379310344827
$ws.Range("M20").Value = -34.615384615384
$ws.Range("N20").Value = -96.551724137931

# ---------------------------------------------------------------------
# Row 21 (G.L.A. - bold TOTAL-style row, styles 17/18)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 12.5
$ws.Range("F21").Value = 44
$ws.Range("G21").Value = 36
$ws.Range("H21").Value = 22.222222222222
$ws.Range("I21").Value = 354
$ws.Range("J21").Value = 346
$ws.Range("K21").Value = 2.312138728323
$ws.Range("L21").Value = -4.838709677419
$ws.Range("M21").Value = 42.741935483871
$ws.Range("N21").Value = -63.877551020408

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 9
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 42
$ws.Range("G24").Value = 21
$ws.Range("H24").Value = 100
$ws.Range("I24").Value = 302
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 20.8
$ws.Range("L24").Value = -8.206686930091
$ws.Range("M24").Value = -20.316622691292

# ---------------------------------------------------------------------
# Row 25 (Retail Theft)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 1
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 137.5
$ws.Range("I25").Value = 178
$ws.Range("J25").Value = 114
$ws.Range("K25").Value = 56.140350877193
$ws.Range("L25").Value = 22.758620689655

# ---------------------------------------------------------------------
# Row 26 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 17
$ws.Range("H26").Value = 30.769230769230
$ws.Range("I26").Value = 143
$ws.Range("J26").Value = 110
$ws.Range("K26").Value = 30
$ws.Range("L26").Value = 1.418439716312
$ws.Range("M26").Value = -4.026845637583

# ---------------------------------------------------------------------
# Row 27 (UCR Rape*): G27 1 -> "0", H27 0 -> "***.*"
# ---------------------------------------------------------------------
Set-TextCell $ws.Range("G27") "0" $donorText0
Set-TextCell $ws.Range("H27") "***.*" $donorTextStar

# ---------------------------------------------------------------------
# Row 28 (Other Sex Crimes): C28 2 -> "0", F28 3 -> 2
# ---------------------------------------------------------------------
Set-TextCell $ws.Range("C28") "0" $donorText0
$ws.Range("F28").Value = 2

# ---------------------------------------------------------------------
# Row 33 (Hate Crimes): D33/E33, G33/H33, J33/K33 text -> numbers
# ---------------------------------------------------------------------
Set-NumberCell $ws.Range("D33") 1 $donorInt
Set-NumberCell $ws.Range("E33") -100 $donorPct
Set-NumberCell $ws.Range("G33") 1 $donorInt
Set-NumberCell $ws.Range("H33") -100 $donorPct
Set-NumberCell $ws.Range("J33") 1 $donorInt
Set-NumberCell $ws.Range("K33") -100 $donorPct

Write-Output "edit complete"
